# Atualizando base de dados da pesquisa via Streamlit
# Append the new survey response row (row 18) to Sheet1, mirroring the
# structure of the previous rows (one column per survey field; columns
# with no answer for this submission are left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

# Columns with no value for this submission - cleared/left blank, same
# as every other "not answered" cell already in the sheet.
$blankCols = @(2,3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,20)
foreach ($col in $blankCols) {
    $ws.Cells.Item($row, $col).Value = ""
}

# A18 - E-mail MRV
$ws.Cells.Item($row, 1).Value = "teste2@mrv.com.br"

# I18 - Data/Hora do Envio
$ws.Cells.Item($row, 9).Value = "2025-05-21 12:09:20"

# U18 - Painéis
$ws.Cells.Item($row, 21).Value = "PAP - Dossiê: Comentário teste 5"

# V18 - Ferramentas
$ws.Cells.Item($row, 22).Value = "Planilha automatizada - teste 2,Objetivo 10,Excel,OUTROS,🟢 Pouco Importante,6.0"
